$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.239.62'
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = '1.829.10'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6124'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -4.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2815'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("E10").Value = '  -5.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07675'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '1.822.67'
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.822'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001010'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6325'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.34%  '
$ws.Range("D16").Value = '2.067.69'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.863'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.94%  '
$ws.Range("D19").Value = '29.233.56'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '227.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.020'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.42%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1313'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.050'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.04%  '
$ws.Range("E28").Value = '  -4.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.480'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06395'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.30%  '
$ws.Range("E31").Value = '  -1.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.827'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.808'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.132'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.747'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6503'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  -1.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.752'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("D39").Value = '1.217.07'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.576'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01741'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9184'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").Value = '1.980.63'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.93'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.61%  '
$ws.Range("E47").Value = '  -1.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.622'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.612'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4570'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05526'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.53%  '
